$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 3 (shifts existing rows 3..21 down to 4..22;
# the last row falls outside the table's original extent and is dropped)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the new IPO entry
$ws.Cells.Item(3, 1).Value = "티디에스팜"
$ws.Cells.Item(3, 2).Value = "2024.07.10~07.16"
$ws.Cells.Item(3, 3).Value = "9,500~10,700"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = 9500
$ws.Cells.Item(3, 6).Value = "한국투자증권"

# Remove the last data row (row 22) so the table keeps 20 data rows (rows 2-21)
$ws.Rows.Item(22).Delete()

# Update the "하스" row's 확정공모가(D) value from "-" to "16000"
$ws.Cells.Item(13, 4).Value = "16000"
